$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Header date
Replace-Text "2023-11-05 Sunday" "2023-11-06 Monday"

# Row 0
Replace-Text "29×96=" "15×15="
Replace-Text "63×57=" "51×72="
Replace-Text "22×73=" "70×48="
Replace-Text "36×96=" "77×16="
Replace-Text "60×98=" "46×11="

# Row 4
Replace-Text "60×32=" "50×35="
Replace-Text "68×75=" "87×74="
Replace-Text "88×85=" "84×32="
Replace-Text "99×81=" "79×80="
Replace-Text "44×40=" "74×18="

# Row 9
Replace-Text "71×89=" "96×49="
Replace-Text "76×91=" "96×39="
Replace-Text "74×72=" "20×48="
Replace-Text "91×62=" "12×75="
Replace-Text "47×47=" "25×70="

# Row 14
Replace-Text "61×53=" "24×24="
Replace-Text "31×45=" "51×15="
Replace-Text "71×60=" "26×28="
Replace-Text "21×68=" "26×65="
Replace-Text "83×91=" "77×40="

# Row 19
Replace-Text "86×89=" "19×95="
Replace-Text "83×99=" "63×12="
Replace-Text "25×90=" "75×51="
Replace-Text "18×95=" "83×83="
Replace-Text "39×95=" "38×40="

Write-Host "Done applying replacements"
